$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Expected Response" (D), add "Actual Response" (E) and "Pass/Fail" (F)
# columns for the chatbot test results table (rows 2-16).

$ws.Range("D2").Value2 = "`"Hello! How can I help you today?`""
$ws.Range("E2").Value2 = "Hello / Hi Good day! What can I do for you today?"
$ws.Range("F2").Value2 = "PASS"
$ws.Range("D3").Value2 = "We are open from 9 AM to 5 PM, Monday to Friday."
$ws.Range("E3").Value2 = "What are your working hours? We are open 9  to 5 , Monday to Friday."
$ws.Range("F3").Value2 = "FAIL"
$ws.Range("D4").Value2 = "We are located at [Company Address]."
$ws.Range("E4").Value2 = "Where are you located? Sorry, could you say that again?"
$ws.Range("F4").Value2 = "PASS"
$ws.Range("D5").Value2 = "You can reach us at +94-XXXXXXXXX"
$ws.Range("E5").Value2 = "What is your contact number? You can email support@example.com or call 123-456-7890."
$ws.Range("F5").Value2 = "PASS"
$ws.Range("D6").Value2 = "You're welcome!"
$ws.Range("E6").Value2 = "Thanks / Thank you I didn't get that. Can you repeat?"
$ws.Range("F6").Value2 = "PASS"
$ws.Range("D7").Value2 = "I'm sorry, I didn’t understand that. Can you rephrase?"
$ws.Range("E7").Value2 = "Blah blah blah What was that?"
$ws.Range("F7").Value2 = "FAIL"
$ws.Range("D8").Value2 = "We are open from 9 AM to 5 PM, Monday to Friday."
$ws.Range("E8").Value2 = "Wht are ur hours? We are open 9  to 5 , Monday to Friday."
$ws.Range("F8").Value2 = "FAIL"
$ws.Range("D9").Value2 = "We offer [Service A], [Service B], and [Service C]."
$ws.Range("E9").Value2 = "What services do you offer? Can you say that again?"
$ws.Range("F9").Value2 = "PASS"
$ws.Range("D10").Value2 = "The price of [Service] starts at `$XX"
$ws.Range("E10").Value2 = "What is the price of [Service]? One more time?"
$ws.Range("F10").Value2 = "FAIL"
$ws.Range("D11").Value2 = "You can contact support at support@email.com"
$ws.Range("E11").Value2 = "How can I contact support? For support, email support@example.com or call 123-456-7890."
$ws.Range("F11").Value2 = "PASS"
$ws.Range("D12").Value2 = "You can email us at info@email.com`n."
$ws.Range("E12").Value2 = "What is your email? Sorry, what was that?"
$ws.Range("F12").Value2 = "PASS"
$ws.Range("D13").Value2 = "Hello! How can I help you today?"
$ws.Range("E13").Value2 = "Hey there! Hello! How can I help you today?"
$ws.Range("F13").Value2 = "PASS"
$ws.Range("D14").Value2 = "We close at 5 PM, Monday to Friday."
$ws.Range("E14").Value2 = "When do you close? Sorry, could you say that again?"
$ws.Range("F14").Value2 = "PASS"
$ws.Range("D15").Value2 = "We open at 9 AM, Monday to Friday."
$ws.Range("E15").Value2 = "When do you open? We are open 9  to 5 , Monday to Friday."
$ws.Range("F15").Value2 = "FAIL"
$ws.Range("D16").Value2 = "Can you please type your question?"
$ws.Range("E16").Value2 = "`"`" (no input) I didn't get that. Can you say it again?"
$ws.Range("F16").Value2 = "PASS"

# Row 12 holds a multi-line actual response; ensure wrap text formatting
# (matches the existing row height already set to accommodate two lines).
$ws.Range("D12").WrapText = $true

Write-Host "Chatbot test results updated."
